$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 36
$ws.Range("C2").Value = 103
$ws.Range("D2").Value = 71
$ws.Range("E2").Value = 8

$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 1

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = 6

$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 40
$ws.Range("D7").Value = 73
$ws.Range("E7").Value = 27

$ws.Range("C8").Value = 74

$ws.Range("C9").Value = 18
$ws.Range("D9").Value = 6

$ws.Range("D10").Value = 93
$ws.Range("E10").Value = 9

$ws.Range("B11").Value = 36
$ws.Range("C11").Value = 137
$ws.Range("D11").Value = 147
$ws.Range("E11").Value = 38

$ws.Range("D12").Value = 75
$ws.Range("E12").Value = 2

$ws.Range("C13").Value = 53
$ws.Range("D13").Value = 37
$ws.Range("E13").Value = 2

$ws.Range("D15").Value = 34

$ws.Range("C17").Value = 47

$wb.Save()
